$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume table (columns D=Price, E=Volume(1h))
# for each ranked row 2-51 per the latest scrape refresh. Row 47/48
# (Cosmos/Filecoin) swapped rank order, so B/C (Coin/Link) are also updated
# for those two rows.
#
# Column D values are plain text (e.g. "96.692.47", "1.00", "0.0000251")
# rather than numbers, so they're entered with a leading apostrophe to force
# Excel to store them as text and keep the exact original formatting/digits
# instead of normalizing them into numeric values.

$ws.Range('D2').Value = "'96.692.47"
$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('D3').Value = "'3.331.49"
$ws.Range('E3').Value = '  -2.41%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'248.24"
$ws.Range('E5').Value = '  -3.21%  '
$ws.Range('D6').Value = "'652.16"
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('E7').Value = '  -7.28%  '
$ws.Range('D8').Value = "'0.418"
$ws.Range('E8').Value = '  -3.42%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = "'0.986"
$ws.Range('E10').Value = '  -7.19%  '
$ws.Range('D11').Value = "'3.330.89"
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('E12').Value = '  -3.87%  '
$ws.Range('D13').Value = "'39.67"
$ws.Range('E13').Value = '  -6.57%  '
$ws.Range('E14').Value = '  -4.88%  '
$ws.Range('D15').Value = "'96.312.63"
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = "'0.0000251"
$ws.Range('E16').Value = '  -3.36%  '
$ws.Range('D17').Value = "'3.950.86"
$ws.Range('E17').Value = '  -2.09%  '
$ws.Range('D18').Value = "'8.46"
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').Value = "'3.335.43"
$ws.Range('E19').Value = '  -2.31%  '
$ws.Range('D20').Value = "'16.75"
$ws.Range('E20').Value = '  -5.01%  '
$ws.Range('D21').Value = "'0.482"
$ws.Range('E21').Value = '  -7.35%  '
$ws.Range('D22').Value = "'504.77"
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').Value = "'10.44"
$ws.Range('E23').Value = '  -4.24%  '
$ws.Range('D24').Value = "'3.34"
$ws.Range('E24').Value = '  -3.72%  '
$ws.Range('E25').Value = '  -4.75%  '
$ws.Range('E26').Value = '  +4.14%  '
$ws.Range('D27').Value = "'94.61"
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('D28').Value = "'11.94"
$ws.Range('E28').Value = '  -6.23%  '
$ws.Range('D29').Value = "'3.496.10"
$ws.Range('E29').Value = '  -1.91%  '
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('E31').Value = '  -7.16%  '
$ws.Range('E32').Value = '  -5.61%  '
$ws.Range('E33').Value = '  -4.83%  '
$ws.Range('E34').Value = '  +12.58%  '
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('D36').Value = "'0.542"
$ws.Range('E36').Value = '  -5.32%  '
$ws.Range('D37').Value = "'27.97"
$ws.Range('E37').Value = '  -6.14%  '
$ws.Range('E38').Value = '  +8.09%  '
$ws.Range('D39').Value = "'7.51"
$ws.Range('E39').Value = '  -4.02%  '
$ws.Range('E41').Value = '  -4.64%  '
$ws.Range('D42').Value = "'505.73"
$ws.Range('E42').Value = '  -1.69%  '
$ws.Range('D43').Value = "'24.50"
$ws.Range('E43').Value = '  -0.86%  '
$ws.Range('D44').Value = "'3.65"
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').Value = "'0.824"
$ws.Range('E45').Value = '  -3.16%  '
$ws.Range('D46').Value = "'0.0405"
$ws.Range('E46').Value = '  -6.76%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = "'5.43"
$ws.Range('E47').Value = '  +1.21%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').Value = "'8.32"
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').Value = "'1.61"
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').Value = "'52.99"
$ws.Range('E50').Value = '  +4.82%  '
$ws.Range('E51').Value = '  -6.06%  '
